$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-11 Friday" "2025-07-12 Saturday"

Replace-Text "41×31=" "14×15="
Replace-Text "73×68=" "55×72="
Replace-Text "19×49=" "55×60="
Replace-Text "36×42=" "26×76="
Replace-Text "32×23=" "51×37="
Replace-Text "15×20=" "27×80="
Replace-Text "17×32=" "74×95="
Replace-Text "63×14=" "40×12="
Replace-Text "17×52=" "95×34="
Replace-Text "15×24=" "84×91="
Replace-Text "48×36=" "15×56="
Replace-Text "24×37=" "62×26="
Replace-Text "38×84=" "24×33="
Replace-Text "38×40=" "67×24="
Replace-Text "50×42=" "43×16="
Replace-Text "65×54=" "22×83="
Replace-Text "87×19=" "90×77="
Replace-Text "70×73=" "77×43="
Replace-Text "86×15=" "60×55="
Replace-Text "96×99=" "65×74="
Replace-Text "92×55=" "21×72="
Replace-Text "15×23=" "86×90="
Replace-Text "40×34=" "17×50="
Replace-Text "75×52=" "58×34="
Replace-Text "14×25=" "49×19="
